# =====================================================================
# Edit script for web-scraping/TP2/test-quotes.xlsx
#
# Commit intent: 'add better logging + scraping page1 + add README.md +
# transforms information in logging (top 5 author, mean citations, top 10 tags)'
#
# The scraper now also pulls page 1 of quotes.toscrape.com (10 more quotes,
# their authors and tags) and prepends them to the existing dataset, which
# pushes the previous rows down and truncates the tail so every sheet keeps
# a fixed row budget. A new 'longueur moyenne' column is added next to the
# citations with the computed mean citation length.
# =====================================================================

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet 'Citations': column A (quotes) + new column B (mean length)
# ---------------------------------------------------------------
$wsCitations = $wb.Worksheets.Item("Citations")

$wsCitations.Range("A1").Value = 'Citations'
$wsCitations.Range("A2").Value = '“The world as we have created it is a process of our thinking. It cannot be changed without changing our thinking.”'
$wsCitations.Range("A3").Value = '“It is our choices, Harry, that show what we truly are, far more than our abilities.”'
$wsCitations.Range("A4").Value = '“There are only two ways to live your life. One is as though nothing is a miracle. The other is as though everything is a miracle.”'
$wsCitations.Range("A5").Value = '“The person, be it gentleman or lady, who has not pleasure in a good novel, must be intolerably stupid.”'
$wsCitations.Range("A6").Value = '“Imperfection is beauty, madness is genius and it''s better to be absolutely ridiculous than absolutely boring.”'
$wsCitations.Range("A7").Value = '“Try not to become a man of success. Rather become a man of value.”'
$wsCitations.Range("A8").Value = '“It is better to be hated for what you are than to be loved for what you are not.”'
$wsCitations.Range("A9").Value = '“I have not failed. I''ve just found 10,000 ways that won''t work.”'
$wsCitations.Range("A10").Value = '“A woman is like a tea bag; you never know how strong it is until it''s in hot water.”'
$wsCitations.Range("A11").Value = '“A day without sunshine is like, you know, night.”'
$wsCitations.Range("A12").Value = '“This life is what you make it. No matter what, you''re going to mess up sometimes, it''s a universal truth. But the good part is you get to decide how you''re going to mess it up. Girls will be your friends - they''ll act like it anyway. But just remember, some come, some go. The ones that stay with you through everything - they''re your true best friends. Don''t let go of them. Also remember, sisters make the best friends in the world. As for lovers, well, they''ll come and go too. And baby, I hate to say it, most of them - actually pretty much all of them are going to break your heart, but you can''t give up because if you give up, you''ll never find your soulmate. You''ll never find that half who makes you whole and that goes for everything. Just because you fail once, doesn''t mean you''re gonna fail at everything. Keep trying, hold on, and always, always, always believe in yourself, because if you don''t, then who will, sweetie? So keep your head high, keep your chin up, and most importantly, keep smiling, because life''s a beautiful thing and there''s so much to smile about.”'
$wsCitations.Range("A13").Value = '“It takes a great deal of bravery to stand up to our enemies, but just as much to stand up to our friends.”'
$wsCitations.Range("A14").Value = '“If you can''t explain it to a six year old, you don''t understand it yourself.”'
$wsCitations.Range("A15").Value = '“You may not be her first, her last, or her only. She loved before she may love again. But if she loves you now, what else matters? She''s not perfect—you aren''t either, and the two of you may never be perfect together but if she can make you laugh, cause you to think twice, and admit to being human and making mistakes, hold onto her and give her the most you can. She may not be thinking about you every second of the day, but she will give you a part of her that she knows you can break—her heart. So don''t hurt her, don''t change her, don''t analyze and don''t expect more than she can give. Smile when she makes you happy, let her know when she makes you mad, and miss her when she''s not there.”'
$wsCitations.Range("A16").Value = '“I like nonsense, it wakes up the brain cells. Fantasy is a necessary ingredient in living.”'
$wsCitations.Range("A17").Value = '“I may not have gone where I intended to go, but I think I have ended up where I needed to be.”'
$wsCitations.Range("A18").Value = '“The opposite of love is not hate, it''s indifference. The opposite of art is not ugliness, it''s indifference. The opposite of faith is not heresy, it''s indifference. And the opposite of life is not death, it''s indifference.”'
$wsCitations.Range("A19").Value = '“It is not a lack of love, but a lack of friendship that makes unhappy marriages.”'
$wsCitations.Range("A20").Value = '“Good friends, good books, and a sleepy conscience: this is the ideal life.”'
$wsCitations.Range("A21").Value = '“Life is what happens to us while we are making other plans.”'
$wsCitations.Range("A22").Value = '“I love you without knowing how, or when, or from where. I love you simply, without problems or pride: I love you in this way because I do not know any other way of loving but this, in which there is no I or you, so intimate that your hand upon my chest is my hand, so intimate that when I fall asleep your eyes close.”'
$wsCitations.Range("A23").Value = '“For every minute you are angry you lose sixty seconds of happiness.”'
$wsCitations.Range("A24").Value = '“If you judge people, you have no time to love them.”'
$wsCitations.Range("A25").Value = '“Anyone who thinks sitting in church can make you a Christian must also think that sitting in a garage can make you a car.”'
$wsCitations.Range("A26").Value = '“Beauty is in the eye of the beholder and it may be necessary from time to time to give a stupid or misinformed beholder a black eye.”'
$wsCitations.Range("A27").Value = '“Today you are You, that is truer than true. There is no one alive who is Youer than You.”'
$wsCitations.Range("A28").Value = '“If you want your children to be intelligent, read them fairy tales. If you want them to be more intelligent, read them more fairy tales.”'
$wsCitations.Range("A29").Value = '“It is impossible to live without failing at something, unless you live so cautiously that you might as well not have lived at all - in which case, you fail by default.”'
$wsCitations.Range("A30").Value = '“Logic will get you from A to Z; imagination will get you everywhere.”'
$wsCitations.Range("A31").Value = '“One good thing about music, when it hits you, you feel no pain.”'

$wsCitations.Range("B1").Value = 'longueur moyenne'
# match B1's style to A1's header style (bold, centered, bordered)
$wsCitations.Range("A1").Copy()
$wsCitations.Range("B1").PasteSpecial(-4122)
$wsCitations.Range("B2").Value = 157.3666666666667
$wsCitations.Range("B3").Value = 157.3666666666667
$wsCitations.Range("B4").Value = 157.3666666666667
$wsCitations.Range("B5").Value = 157.3666666666667
$wsCitations.Range("B6").Value = 157.3666666666667
$wsCitations.Range("B7").Value = 157.3666666666667
$wsCitations.Range("B8").Value = 157.3666666666667
$wsCitations.Range("B9").Value = 157.3666666666667
$wsCitations.Range("B10").Value = 157.3666666666667
$wsCitations.Range("B11").Value = 157.3666666666667
$wsCitations.Range("B12").Value = 157.3666666666667
$wsCitations.Range("B13").Value = 157.3666666666667
$wsCitations.Range("B14").Value = 157.3666666666667
$wsCitations.Range("B15").Value = 157.3666666666667
$wsCitations.Range("B16").Value = 157.3666666666667
$wsCitations.Range("B17").Value = 157.3666666666667
$wsCitations.Range("B18").Value = 157.3666666666667
$wsCitations.Range("B19").Value = 157.3666666666667
$wsCitations.Range("B20").Value = 157.3666666666667
$wsCitations.Range("B21").Value = 157.3666666666667
$wsCitations.Range("B22").Value = 157.3666666666667
$wsCitations.Range("B23").Value = 157.3666666666667
$wsCitations.Range("B24").Value = 157.3666666666667
$wsCitations.Range("B25").Value = 157.3666666666667
$wsCitations.Range("B26").Value = 157.3666666666667
$wsCitations.Range("B27").Value = 157.3666666666667
$wsCitations.Range("B28").Value = 157.3666666666667
$wsCitations.Range("B29").Value = 157.3666666666667
$wsCitations.Range("B30").Value = 157.3666666666667
$wsCitations.Range("B31").Value = 157.3666666666667

# ---------------------------------------------------------------
# Sheet 'Auteurs': column A (authors, aligned with new citations)
# ---------------------------------------------------------------
$wsAuteurs = $wb.Worksheets.Item("Auteurs")

$wsAuteurs.Range("A1").Value = 'Auteurs'
$wsAuteurs.Range("A2").Value = 'Albert Einstein'
$wsAuteurs.Range("A3").Value = 'J.K. Rowling'
$wsAuteurs.Range("A4").Value = 'Albert Einstein'
$wsAuteurs.Range("A5").Value = 'Jane Austen'
$wsAuteurs.Range("A6").Value = 'Marilyn Monroe'
$wsAuteurs.Range("A7").Value = 'Albert Einstein'
$wsAuteurs.Range("A8").Value = 'André Gide'
$wsAuteurs.Range("A9").Value = 'Thomas A. Edison'
$wsAuteurs.Range("A10").Value = 'Eleanor Roosevelt'
$wsAuteurs.Range("A11").Value = 'Steve Martin'
$wsAuteurs.Range("A12").Value = 'Marilyn Monroe'
$wsAuteurs.Range("A13").Value = 'J.K. Rowling'
$wsAuteurs.Range("A14").Value = 'Albert Einstein'
$wsAuteurs.Range("A15").Value = 'Bob Marley'
$wsAuteurs.Range("A16").Value = 'Dr. Seuss'
$wsAuteurs.Range("A17").Value = 'Douglas Adams'
$wsAuteurs.Range("A18").Value = 'Elie Wiesel'
$wsAuteurs.Range("A19").Value = 'Friedrich Nietzsche'
$wsAuteurs.Range("A20").Value = 'Mark Twain'
$wsAuteurs.Range("A21").Value = 'Allen Saunders'
$wsAuteurs.Range("A22").Value = 'Pablo Neruda'
$wsAuteurs.Range("A23").Value = 'Ralph Waldo Emerson'
$wsAuteurs.Range("A24").Value = 'Mother Teresa'
$wsAuteurs.Range("A25").Value = 'Garrison Keillor'
$wsAuteurs.Range("A26").Value = 'Jim Henson'
$wsAuteurs.Range("A27").Value = 'Dr. Seuss'
$wsAuteurs.Range("A28").Value = 'Albert Einstein'
$wsAuteurs.Range("A29").Value = 'J.K. Rowling'
$wsAuteurs.Range("A30").Value = 'Albert Einstein'
$wsAuteurs.Range("A31").Value = 'Bob Marley'

# ---------------------------------------------------------------
# Sheet 'Tags': column A (tags, aligned with new citations)
# ---------------------------------------------------------------
$wsTags = $wb.Worksheets.Item("Tags")

$wsTags.Range("A1").Value = 'Tags'
$wsTags.Range("A2").Value = 'change'
$wsTags.Range("A3").Value = 'deep-thoughts'
$wsTags.Range("A4").Value = 'thinking'
$wsTags.Range("A5").Value = 'world'
$wsTags.Range("A6").Value = 'abilities'
$wsTags.Range("A7").Value = 'choices'
$wsTags.Range("A8").Value = 'inspirational'
$wsTags.Range("A9").Value = 'life'
$wsTags.Range("A10").Value = 'live'
$wsTags.Range("A11").Value = 'miracle'
$wsTags.Range("A12").Value = 'miracles'
$wsTags.Range("A13").Value = 'aliteracy'
$wsTags.Range("A14").Value = 'books'
$wsTags.Range("A15").Value = 'classic'
$wsTags.Range("A16").Value = 'humor'
$wsTags.Range("A17").Value = 'be-yourself'
$wsTags.Range("A18").Value = 'inspirational'
$wsTags.Range("A19").Value = 'adulthood'
$wsTags.Range("A20").Value = 'success'
$wsTags.Range("A21").Value = 'value'
$wsTags.Range("A22").Value = 'life'
$wsTags.Range("A23").Value = 'love'
$wsTags.Range("A24").Value = 'edison'
$wsTags.Range("A25").Value = 'failure'
$wsTags.Range("A26").Value = 'inspirational'
$wsTags.Range("A27").Value = 'paraphrased'
$wsTags.Range("A28").Value = 'misattributed-eleanor-roosevelt'
$wsTags.Range("A29").Value = 'humor'
$wsTags.Range("A30").Value = 'obvious'
$wsTags.Range("A31").Value = 'simile'
$wsTags.Range("A32").Value = 'friends'
$wsTags.Range("A33").Value = 'heartbreak'
$wsTags.Range("A34").Value = 'inspirational'
$wsTags.Range("A35").Value = 'life'
$wsTags.Range("A36").Value = 'love'
$wsTags.Range("A37").Value = 'sisters'
$wsTags.Range("A38").Value = 'courage'
$wsTags.Range("A39").Value = 'friends'
$wsTags.Range("A40").Value = 'simplicity'
$wsTags.Range("A41").Value = 'understand'
$wsTags.Range("A42").Value = 'love'
$wsTags.Range("A43").Value = 'fantasy'
$wsTags.Range("A44").Value = 'life'
$wsTags.Range("A45").Value = 'navigation'
$wsTags.Range("A46").Value = 'activism'
$wsTags.Range("A47").Value = 'apathy'
$wsTags.Range("A48").Value = 'hate'
$wsTags.Range("A49").Value = 'indifference'
$wsTags.Range("A50").Value = 'inspirational'
$wsTags.Range("A51").Value = 'love'
$wsTags.Range("A52").Value = 'opposite'
$wsTags.Range("A53").Value = 'philosophy'
$wsTags.Range("A54").Value = 'friendship'
$wsTags.Range("A55").Value = 'lack-of-friendship'
$wsTags.Range("A56").Value = 'lack-of-love'
$wsTags.Range("A57").Value = 'love'
$wsTags.Range("A58").Value = 'marriage'
$wsTags.Range("A59").Value = 'unhappy-marriage'
$wsTags.Range("A60").Value = 'books'
$wsTags.Range("A61").Value = 'contentment'
$wsTags.Range("A62").Value = 'friends'
$wsTags.Range("A63").Value = 'friendship'
$wsTags.Range("A64").Value = 'life'
$wsTags.Range("A65").Value = 'fate'
$wsTags.Range("A66").Value = 'life'
$wsTags.Range("A67").Value = 'misattributed-john-lennon'
$wsTags.Range("A68").Value = 'planning'
$wsTags.Range("A69").Value = 'plans'
$wsTags.Range("A70").Value = 'love'
$wsTags.Range("A71").Value = 'poetry'
$wsTags.Range("A72").Value = 'happiness'
$wsTags.Range("A73").Value = 'attributed-no-source'
$wsTags.Range("A74").Value = 'humor'
$wsTags.Range("A75").Value = 'religion'
$wsTags.Range("A76").Value = 'humor'
$wsTags.Range("A77").Value = 'comedy'
$wsTags.Range("A78").Value = 'life'
$wsTags.Range("A79").Value = 'yourself'
$wsTags.Range("A80").Value = 'children'
$wsTags.Range("A81").Value = 'fairy-tales'
$wsTags.Range("A82").Value = 'imagination'
$wsTags.Range("A83").Value = 'music'

